$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B24").Value = -0.99949039209519
$ws.Range("C24").Value = -167966135.730231
$ws.Range("D24").Value = 0.95662930260746
$ws.Range("E24").Value = 0.999490392095189
$ws.Range("F24").Value = -0.95662930260746
$ws.Range("G24").Value = 123411.461836652
$ws.Range("H24").Value = 7054577.7006697
$ws.Range("I24").Value = -2938368.1389679
$ws.Range("J24").Value = 0.956629302607457
$ws.Range("K24").Value = 0.999490392095189
$ws.Range("L24").Value = -0.95662930260746
$ws.Range("M24").Value = 50.2489665458678
$ws.Range("N24").Value = 2872.38505727593
$ws.Range("O24").Value = -2938368.1389679
